# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (Stock) worksheet (5th tab) gets three new trailing columns:
#   H = date, I = legislator_name, J = legislator_id
# with a matching data row for the single stock entry already present.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# ---- Header row (row 1): copy the existing bold/bordered/centered header
# style from column B onto the three new header cells, then set their text.
$ws.Cells.Item(1, 2).Copy()
$ws.Range($ws.Cells.Item(1, 8), $ws.Cells.Item(1, 10)).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

# ---- Data row (row 2): fill in the values for the existing stock entry.
# Force the date cell to be stored as text (not auto-parsed into a date
# serial number) by giving it a text number format before assigning it.
$ws.Cells.Item(2, 8).NumberFormat = "@"
$ws.Cells.Item(2, 8).Value = "2012-04-27"
$ws.Cells.Item(2, 9).Value = "許智傑"
$ws.Cells.Item(2, 10).Value = 1750
